$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "PAZ ANASTACIO JUANITA ROSA",
    "ARRUNATEGUI ESPINOZA JOVANNY",
    "CARRION LAZARO MICHAEL LUIS",
    "NIMA CARMEN KAREN DEL MILAGRO",
    "ALZAMORA CHERRES SIRLEY YASMIN",
    "ESPINOZA VALDIVIEZO JUNIOR RICARDO",
    "NAVARRO JUAREZ LIDIA",
    "PULACHE LAZO VILMA YOHANA",
    "DOMINGUEZ CUEVA MERLING DEL JESUS YOLINDA",
    "LILIAN ROXANA VEGA GARCÍA",
    "NIMA CRUZ ANA GRACIELA"
)

$totals = @(105, 104, 102, 90, 87, 84, 80, 79, 70, 70, 1)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}
